$d = $word.ActiveDocument

$needle = "}}, City of {{sellers_current_city}}, State of {{sellers_current_state"
$text = $d.Content.Text
$idx = $text.IndexOf($needle)

if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + $needle.Length)
    $r.Delete()
}
